$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("subject")

$ws.Rows("3:4").Insert()
$ws.Range("B3").Value = 'do_extract_trigger_duration'
$ws.Range("C3").Value = 'extract /identify triggers from one channel using deflection durations'
$ws.Rows(3).RowHeight = 30
$ws.Range("B4").Value = 'do_extract_trigger_evaluate'
$ws.Range("C4").Value = 'extract /identify triggers from one channel (using its deflections) which are bounded between couples of events and evaluate delayes between events and extracted trigger onsets'
$ws.Rows(4).RowHeight = 60

$ws.Rows("20:21").Insert()
$ws.Range("B20").Value = 'do_recover_asr'
$ws.Range("C20").Value = 'recover asr cleaned file'
$ws.Range("B21").Value = 'do_cleanline'
$ws.Range("C21").Value = 'cleanline'

$ws.Rows("34:34").Insert()
$ws.Range("B34").Value = 'do_compute_hr'
$ws.Range("C34").Value = 'compute heart rate from EKG channel and add it as a channel'
$ws.Rows(34).RowHeight = 30

$ws.Rows("43:44").Insert()
$ws.Range("B43").Value = 'do_darbeliai_export2ragu'
$ws.Range("C43").Value = 'export data from eeglab to ragu'
$ws.Range("B44").Value = 'do_ragu'
$ws.Range("C44").Value = 'process data with ragu'

$ws.Rows("48:51").Insert()
$ws.Range("B48").Value = 'do_export_data'
$ws.Range("C48").Value = 'export data in standard formats like EDF'
$ws.Range("B49").Value = 'do_subject_erp_curve'
$ws.Range("C49").Value = 'plot comparison of erp curve between conditions for single subject'
$ws.Rows(49).RowHeight = 30
$ws.Range("B50").Value = 'do_subject_erp_topo'
$ws.Range("C50").Value = 'plot comparison of erp topographic map between conditions for single subject'
$ws.Rows(50).RowHeight = 30
$ws.Range("B51").Value = 'do_subject_ersp_tf'
$ws.Range("C51").Value = 'plot comparison of ersp time-frequency representation between conditions for single subject'
$ws.Rows(51).RowHeight = 45

# Fix up column A auto-increment formulas for the whole range (insert can leave stale refs)
for ($r = 2; $r -le 51; $r++) {
    $prev = $r - 1
    $ws.Range("A$r").Formula = "=A$prev+1"
}

$excel.CalculateFull()
$ws.Range("C5").Select()
